$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# User typed "Change A1" into cell A1 (as reflected in the shared-strings
# table / sheet1.xml in the target commit) and then moved on to A2, which
# is exactly what Excel does after committing a cell edit with Enter.
$ws.Range("A1").Value = "Change A1"
[void]$ws.Range("A2").Select()
